$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.6495703901946965
$ws.Range("D2").Value = 0.5336804966524851
$ws.Range("E2").Value = 0.6598091001760726
$ws.Range("F2").Value = 0.5550511827678311

$ws.Range("C3").Value = 0.6747805441581174
$ws.Range("D3").Value = 0.6277410221681932
$ws.Range("E3").Value = 0.6805671392827356
$ws.Range("F3").Value = 0.6452157957942443

$ws.Range("C4").Value = 0.7085126875091203
$ws.Range("D4").Value = 0.6145867026921692
$ws.Range("E4").Value = 0.7127235659345751
$ws.Range("F4").Value = 0.6385898432540839

$ws.Range("C5").Value = 0.7977292576506659
$ws.Range("D5").Value = 0.7590446704583426
$ws.Range("E5").Value = 0.7795385043091465
$ws.Range("F5").Value = 0.7687835241052147

$ws.Range("C6").Value = 0.8133977704366796
$ws.Range("D6").Value = 0.7876622016340096
$ws.Range("E6").Value = 0.8059494022796775
$ws.Range("F6").Value = 0.7992097196385826

$ws.Range("C7").Value = 0.8045213424983151
$ws.Range("D7").Value = 0.7636442298603677
$ws.Range("E7").Value = 0.7972384394402743
$ws.Range("F7").Value = 0.7817481003397651

$ws.Range("C8").Value = 0.7717053046795664
$ws.Range("D8").Value = 0.7320382797320696
$ws.Range("E8").Value = 0.7740709850801594
$ws.Range("F8").Value = 0.7494307527014981
